$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (Company Name), B (Company Number), H (Category)
# for data rows 2-8, reflecting the reordering described in the diff.
$data = @(
    @{ Row = 2;  A = "T GILPIN PHYSIO CONSULTANCY LTD";  B = "16460503";  H = "LP" },
    @{ Row = 3;  A = "SAMVIV PARTNERS LTD";               B = "16460672";  H = "Partners" },
    @{ Row = 4;  A = "4D CAPITAL PROPCO (44) LIMITED";    B = "16461269";  H = "Capital" },
    @{ Row = 5;  A = "DGPI LTD";                          B = "SC849118";  H = "GP" },
    @{ Row = 6;  A = "DAVIDSON CAPITAL HOLDINGS LTD";     B = "SC849117";  H = "Capital" },
    @{ Row = 7;  A = "AFROSCOT VENTURES LTD";             B = "16462878";  H = "Ventures" },
    @{ Row = 8;  A = "ST GEORGE CAPITAL (LAND) LIMITED";  B = "16462880";  H = "Capital" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.A

    # Force the company-number column to be stored as text (matching the
    # original inline-string typing) instead of being auto-converted to a
    # number, then clear the formatting tweak so no stray style is left
    # behind on the cell.
    $ws.Range("B$r").NumberFormat = "@"
    $ws.Range("B$r").Value = $item.B
    $ws.Range("B$r").ClearFormats()

    $ws.Range("H$r").Value = $item.H
}
